$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update financial figures
$ws.Range("D2").Value = 12694
$ws.Range("E2").Value = 246
$ws.Range("F2").Value = 231
$ws.Range("G2").Value = 438
$ws.Range("H2").Value = 270
$ws.Range("I2").Value = 178
$ws.Range("J2").Value = 92
$ws.Range("K2").Value = 11902
$ws.Range("L2").Value = 6717
$ws.Range("M2").Value = 5185
$ws.Range("N2").Value = 3674
$ws.Range("O2").Value = 1511
$ws.Range("P2").Value = 161
$ws.Range("Q2").Value = 440
$ws.Range("R2").Value = -336
$ws.Range("S2").Value = -103
$ws.Range("T2").Value = 364
$ws.Range("U2").Value = 76
$ws.Range("V2").Value = 2763
$ws.Range("W2").Value = 1.94
$ws.Range("X2").Value = 2.13
$ws.Range("Y2").Value = 4.94
$ws.Range("Z2").Value = 2.32
$ws.Range("AA2").Value = 129.53
$ws.Range("AB2").Value = 2318.38
$ws.Range("AC2").Value = 1104
$ws.Range("AD2").Value = 9.69
$ws.Range("AE2").Value = 22836
$ws.Range("AF2").Value = 0.47
$ws.Range("AG2").Value = 250
$ws.Range("AH2").Value = 2.34
$ws.Range("AI2").Value = 22.65
$ws.Range("AJ2").Value = 16089459

# Row 3: update financial figures
$ws.Range("D3").Value = 10560
$ws.Range("E3").Value = 240
$ws.Range("F3").Value = 240
$ws.Range("G3").Value = 339
$ws.Range("H3").Value = 277
$ws.Range("I3").Value = 192
$ws.Range("J3").Value = 84
$ws.Range("K3").Value = 11394
$ws.Range("L3").Value = 6043
$ws.Range("M3").Value = 5351
$ws.Range("N3").Value = 3809
$ws.Range("O3").Value = 1542
$ws.Range("P3").Value = 161
$ws.Range("Q3").Value = 578
$ws.Range("R3").Value = -1128
$ws.Range("S3").Value = -51
$ws.Range("T3").Value = 541
$ws.Range("U3").Value = 37
$ws.Range("V3").Value = 2785
$ws.Range("W3").Value = 2.28
$ws.Range("X3").Value = 2.62
$ws.Range("Y3").Value = 5.15
$ws.Range("Z3").Value = 2.38
$ws.Range("AA3").Value = 112.93
$ws.Range("AB3").Value = 2411.23
$ws.Range("AC3").Value = 1196
$ws.Range("AD3").Value = 8.9
$ws.Range("AE3").Value = 23671
$ws.Range("AF3").Value = 0.45
$ws.Range("AG3").Value = 250
$ws.Range("AH3").Value = 2.35
$ws.Range("AI3").Value = 20.9
$ws.Range("AJ3").Value = 16089459

# Row 4: update financial figures
$ws.Range("D4").Value = 8747
$ws.Range("E4").Value = 269
$ws.Range("F4").Value = 294
$ws.Range("G4").Value = 281
$ws.Range("H4").Value = 234
$ws.Range("I4").Value = 142
$ws.Range("J4").Value = 92
$ws.Range("K4").Value = 11284
$ws.Range("L4").Value = 5842
$ws.Range("M4").Value = 5443
$ws.Range("N4").Value = 3875
$ws.Range("O4").Value = 1568
$ws.Range("P4").Value = 161
$ws.Range("Q4").Value = 629
$ws.Range("R4").Value = -889
$ws.Range("S4").Value = -293
$ws.Range("T4").Value = 544
$ws.Range("U4").Value = 85
$ws.Range("V4").Value = 2600
$ws.Range("W4").Value = 3.08
$ws.Range("X4").Value = 2.67
$ws.Range("Y4").Value = 3.7
$ws.Range("Z4").Value = 2.06
$ws.Range("AA4").Value = 107.33
$ws.Range("AB4").Value = 2461.42
$ws.Range("AC4").Value = 883
$ws.Range("AD4").Value = 11.33
$ws.Range("AE4").Value = 24082
$ws.Range("AF4").Value = 0.42
$ws.Range("AG4").Value = 250
$ws.Range("AH4").Value = 2.5
$ws.Range("AI4").Value = 28.32
$ws.Range("AJ4").Value = 16089459

# Row 5: update financial figures
$ws.Range("D5").Value = 9031
$ws.Range("E5").Value = 245
$ws.Range("F5").Value = 245
$ws.Range("G5").Value = 337
$ws.Range("H5").Value = 322
$ws.Range("I5").Value = 211
$ws.Range("J5").Value = 111
$ws.Range("K5").Value = 11062
$ws.Range("L5").Value = 5851
$ws.Range("M5").Value = 5211
$ws.Range("N5").Value = 4021
$ws.Range("O5").Value = 1190
$ws.Range("P5").Value = 161
$ws.Range("Q5").Value = 634
$ws.Range("R5").Value = 49
$ws.Range("S5").Value = -712
$ws.Range("T5").Value = 531
$ws.Range("U5").Value = 103
$ws.Range("V5").Value = 2418
$ws.Range("W5").Value = 2.72
$ws.Range("X5").Value = 3.56
$ws.Range("Y5").Value = 5.34
$ws.Range("Z5").Value = 2.88
$ws.Range("AA5").Value = 112.27
$ws.Range("AB5").Value = 2559.02
$ws.Range("AC5").Value = 1309
$ws.Range("AD5").Value = 6.06
$ws.Range("AE5").Value = 24992
$ws.Range("AF5").Value = 0.32
$ws.Range("AG5").Value = 250
$ws.Range("AH5").Value = 3.15
$ws.Range("AI5").Value = 19.09
$ws.Range("AJ5").Value = 16089459

# Row 6: update financial figures
$ws.Range("D6").Value = 9098
$ws.Range("E6").Value = 184
$ws.Range("F6").Value = 184
$ws.Range("G6").Value = 301
$ws.Range("H6").Value = 238
$ws.Range("I6").Value = 195
$ws.Range("K6").Value = 11173
$ws.Range("L6").Value = 5847
$ws.Range("M6").Value = 5326
$ws.Range("N6").Value = 4129
$ws.Range("P6").Value = 161
$ws.Range("Q6").Value = 564
$ws.Range("R6").Value = -356
$ws.Range("S6").Value = -169
$ws.Range("T6").Value = 485
$ws.Range("U6").Value = 79
$ws.Range("V6").Value = 2321
$ws.Range("W6").Value = 2.03
$ws.Range("X6").Value = 2.61
$ws.Range("Y6").Value = 4.79
$ws.Range("Z6").Value = 2.14
$ws.Range("AA6").Value = 109.78
$ws.Range("AB6").Value = 2635.72
$ws.Range("AC6").Value = 1214
$ws.Range("AD6").Value = 5.24
$ws.Range("AE6").Value = 25660
$ws.Range("AF6").Value = 0.25
$ws.Range("AG6").Value = 250
$ws.Range("AH6").Value = 3.93
$ws.Range("AI6").Value = 20.59
$ws.Range("AJ6").Value = 16089459

# Row 7: clear stale data cells, keep only A7:C7
$ws.Range("D7").ClearContents()
$ws.Range("E7").ClearContents()
$ws.Range("G7").ClearContents()
$ws.Range("H7").ClearContents()
$ws.Range("I7").ClearContents()
$ws.Range("K7").ClearContents()
$ws.Range("L7").ClearContents()
$ws.Range("M7").ClearContents()
$ws.Range("N7").ClearContents()
$ws.Range("P7").ClearContents()
$ws.Range("Q7").ClearContents()
$ws.Range("R7").ClearContents()
$ws.Range("S7").ClearContents()
$ws.Range("T7").ClearContents()
$ws.Range("U7").ClearContents()
$ws.Range("W7").ClearContents()
$ws.Range("X7").ClearContents()
$ws.Range("Y7").ClearContents()
$ws.Range("Z7").ClearContents()
$ws.Range("AA7").ClearContents()
$ws.Range("AC7").ClearContents()
$ws.Range("AD7").ClearContents()
$ws.Range("AE7").ClearContents()
$ws.Range("AF7").ClearContents()
$ws.Range("AG7").ClearContents()
$ws.Range("AH7").ClearContents()
$ws.Range("AI7").ClearContents()

# Row 8: clear stale data cells, keep only A8:C8
$ws.Range("D8").ClearContents()
$ws.Range("E8").ClearContents()
$ws.Range("G8").ClearContents()
$ws.Range("H8").ClearContents()
$ws.Range("I8").ClearContents()
$ws.Range("K8").ClearContents()
$ws.Range("L8").ClearContents()
$ws.Range("M8").ClearContents()
$ws.Range("N8").ClearContents()
$ws.Range("P8").ClearContents()
$ws.Range("Q8").ClearContents()
$ws.Range("R8").ClearContents()
$ws.Range("S8").ClearContents()
$ws.Range("T8").ClearContents()
$ws.Range("U8").ClearContents()
$ws.Range("W8").ClearContents()
$ws.Range("X8").ClearContents()
$ws.Range("Y8").ClearContents()
$ws.Range("Z8").ClearContents()
$ws.Range("AA8").ClearContents()
$ws.Range("AC8").ClearContents()
$ws.Range("AD8").ClearContents()
$ws.Range("AE8").ClearContents()
$ws.Range("AF8").ClearContents()
$ws.Range("AG8").ClearContents()
$ws.Range("AH8").ClearContents()
$ws.Range("AI8").ClearContents()

# Row 9: clear stale data cells, keep only A9:C9
$ws.Range("D9").ClearContents()
$ws.Range("E9").ClearContents()
$ws.Range("G9").ClearContents()
$ws.Range("H9").ClearContents()
$ws.Range("I9").ClearContents()
$ws.Range("K9").ClearContents()
$ws.Range("L9").ClearContents()
$ws.Range("M9").ClearContents()
$ws.Range("N9").ClearContents()
$ws.Range("P9").ClearContents()
$ws.Range("Q9").ClearContents()
$ws.Range("R9").ClearContents()
$ws.Range("S9").ClearContents()
$ws.Range("T9").ClearContents()
$ws.Range("U9").ClearContents()
$ws.Range("W9").ClearContents()
$ws.Range("X9").ClearContents()
$ws.Range("Y9").ClearContents()
$ws.Range("Z9").ClearContents()
$ws.Range("AA9").ClearContents()
$ws.Range("AC9").ClearContents()
$ws.Range("AD9").ClearContents()
$ws.Range("AE9").ClearContents()
$ws.Range("AF9").ClearContents()
$ws.Range("AG9").ClearContents()
$ws.Range("AH9").ClearContents()
$ws.Range("AI9").ClearContents()
